$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source table")

# "Supported since" column (D): version bump 0.1.6* -> 0.2.0
$ws.Range("D3").Value = "0.2.0"
$ws.Range("D22").Value = "0.2.0"
$ws.Range("D40").Value = "0.2.0"
$ws.Range("D41").Value = "0.2.0"
$ws.Range("D43").Value = "0.2.0"
$ws.Range("D44").Value = "0.2.0"
$ws.Range("D62").Value = "0.2.0"
$ws.Range("D63").Value = "0.2.0"

# "Notes" column (G): fix markdown link syntax / version bump
$ws.Range("G31").Value = "Moved to io.crosssections in 0.2.0"
$ws.Range("G59").Value = "Moved to io.rr in 0.2.0"
$ws.Range("G60").Value = "Moved to io.rr in 0.2.0"
$ws.Range("G3").Value = "Critical bugfix for [#127](https://github.com/Deltares/HYDROLIB-core/issues/127)."

# Update the active selection on the Source table sheet
$ws.Activate()
$ws.Range("G4").Select()
